$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = -1452
# Row 52
$ws.Range("H52").Value = 439.8
$ws.Range("J52").Value = 439.8
$ws.Range("L52").Value = 1319.4
$ws.Range("N52").Value = -1639.4
# Row 57
$ws.Range("H57").Value = 57833
$ws.Range("J57").Value = 60999.5
$ws.Range("L57").Value = 182998.5
$ws.Range("N57").Value = -183996.5
# Row 64
$ws.Range("H64").Value = 8306.299999999999
$ws.Range("J64").Value = 9626.571
$ws.Range("L64").Value = 9626.571
$ws.Range("N64").Value = -10122.571
# Row 67
$ws.Range("H67").Value = 8306.299999999999
$ws.Range("J67").Value = 9626.571
$ws.Range("L67").Value = 9626.571
$ws.Range("N67").Value = -11342.571
# Row 99
$ws.Range("H99").Value = 1791.5
$ws.Range("I99").Value = 450.4
$ws.Range("J99").Value = 4026.6667
$ws.Range("K99").Value = 1351.2
$ws.Range("L99").Value = 12080.0001
$ws.Range("M99").Value = 146.8000000000002
$ws.Range("N99").Value = -15076.0001
# Row 114
$ws.Range("H114").Value = 78722
$ws.Range("J114").Value = 78722
$ws.Range("L114").Value = 78722
$ws.Range("N114").Value = -87400
# Row 132
$ws.Range("H132").Value = 10383.538
$ws.Range("I132").Value = 9582.208000000001
$ws.Range("K132").Value = 28746.624
$ws.Range("M132").Value = -26216.624
# Row 137
$ws.Range("H137").Value = 4549.8096
$ws.Range("I137").Value = 1824.6
$ws.Range("K137").Value = 5473.799999999999
$ws.Range("M137").Value = -2923.799999999999
# Row 138
$ws.Range("H138").Value = 2620.1
$ws.Range("I138").Value = 1953.7142
$ws.Range("J138").Value = 4175
$ws.Range("K138").Value = 5861.142599999999
$ws.Range("L138").Value = 12525
$ws.Range("M138").Value = -721.1425999999992
$ws.Range("N138").Value = -22805

$ws = $wb.Worksheets.Item("ARM")
# Row 28
$ws.Range("H28").Value = 2831.8
$ws.Range("I28").Value = 2831.8
$ws.Range("K28").Value = 2831.8
$ws.Range("M28").Value = -2639.8
# Row 41
$ws.Range("H41").Value = 1762.25
$ws.Range("I41").Value = 1762.25
$ws.Range("K41").Value = 1762.25
$ws.Range("M41").Value = -1348.25
# Row 99
$ws.Range("H99").Value = 2831.8
$ws.Range("I99").Value = 2831.8
$ws.Range("K99").Value = 2831.8
$ws.Range("M99").Value = 163.1999999999998

$ws = $wb.Worksheets.Item("BSM")
# Row 75
$ws.Range("H75").Value = 37416.668
$ws.Range("J75").Value = 54750
$ws.Range("L75").Value = 54750
$ws.Range("N75").Value = -56622
# Row 78
$ws.Range("H78").Value = 37416.668
$ws.Range("J78").Value = 54750
$ws.Range("L78").Value = 164250
$ws.Range("N78").Value = -173610
# Row 92
$ws.Range("H92").Value = 22999.6
$ws.Range("J92").Value = 22999.6
$ws.Range("L92").Value = 22999.6
$ws.Range("N92").Value = -27991.6

$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Range("H107").Value = 669.4167
$ws.Range("I107").Value = 646.55554
$ws.Range("J107").Value = 738
$ws.Range("K107").Value = 646.55554
$ws.Range("L107").Value = 738
$ws.Range("M107").Value = 1273.44446
$ws.Range("N107").Value = -4578
# Row 122
$ws.Range("H122").Value = 2447.5833
$ws.Range("I122").Value = 2324
$ws.Range("J122").Value = 2535.8572
$ws.Range("K122").Value = 6972
$ws.Range("L122").Value = 7607.571599999999
$ws.Range("M122").Value = -4522
$ws.Range("N122").Value = -12507.5716
# Row 132
$ws.Range("H132").Value = 2299.8333
$ws.Range("I132").Value = 1847.1765
$ws.Range("K132").Value = 5541.529500000001
$ws.Range("M132").Value = -3011.529500000001
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 43.666668
$ws.Range("I2").Value = 42.5
$ws.Range("K2").Value = 255
$ws.Range("M2").Value = -142
# Row 8
$ws.Range("H8").Value = 379.2857
$ws.Range("I8").Value = 379.2857
$ws.Range("K8").Value = 1137.8571
$ws.Range("M8").Value = -998.8571000000002
# Row 132
$ws.Range("H132").Value = 2819.1155
$ws.Range("I132").Value = 2728.7144
$ws.Range("J132").Value = 2852.4211
$ws.Range("K132").Value = 24558.4296
$ws.Range("L132").Value = 25671.7899
$ws.Range("M132").Value = -22028.4296
$ws.Range("N132").Value = -30731.7899

$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 8300732.5
$ws.Range("I11").Value = 8538616
$ws.Range("J11").Value = 7457327.5
$ws.Range("K11").Value = 8538616
$ws.Range("L11").Value = 7457327.5
$ws.Range("M11").Value = -8538477
$ws.Range("N11").Value = -7457605.5
# Row 137
$ws.Range("H137").Value = 162699.75
$ws.Range("J137").Value = 199833
$ws.Range("L137").Value = 199833
$ws.Range("N137").Value = -210033
# Row 138
$ws.Range("H138").Value = 99000
$ws.Range("I138").Value = 99000
$ws.Range("K138").Value = 99000
$ws.Range("M138").Value = -93860

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 7615
$ws.Range("I46").Value = 11530.75
$ws.Range("J46").Value = 6496.2144
$ws.Range("K46").Value = 11530.75
$ws.Range("L46").Value = 6496.2144
$ws.Range("M46").Value = -11342.75
$ws.Range("N46").Value = -6872.2144

$ws = $wb.Worksheets.Item("WVR")
# Row 51
$ws.Range("H51").Value = 70000
$ws.Range("J51").Value = 70000
$ws.Range("L51").Value = 70000
$ws.Range("N51").Value = -71020
# Row 52
$ws.Range("H52").Value = 37000
$ws.Range("I52").Value = 37000
$ws.Range("K52").Value = 37000
$ws.Range("M52").Value = -36774
# Row 62
$ws.Range("H62").Value = 10363.538
$ws.Range("I62").Value = 6250
$ws.Range("J62").Value = 11111.454
$ws.Range("K62").Value = 6250
$ws.Range("L62").Value = 11111.454
$ws.Range("M62").Value = -5626
$ws.Range("N62").Value = -12359.454
# Row 63
$ws.Range("H63").Value = 25649.6
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 25649.6
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 25649.6
$ws.Range("N63").Value = -26897.6
$ws.Range("M63").ClearContents()
# Row 65
$ws.Range("H65").Value = 10363.538
$ws.Range("I65").Value = 6250
$ws.Range("J65").Value = 11111.454
$ws.Range("K65").Value = 31250
$ws.Range("L65").Value = 55557.27
$ws.Range("M65").Value = -28130
$ws.Range("N65").Value = -61797.27
# Row 66
$ws.Range("H66").Value = 25649.6
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 25649.6
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 76948.79999999999
$ws.Range("N66").Value = -83188.79999999999
$ws.Range("M66").ClearContents()
# Row 68
$ws.Range("H68").Value = 31949
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 31949
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 31949
$ws.Range("N68").Value = -33571
$ws.Range("M68").ClearContents()
# Row 71
$ws.Range("H71").Value = 31949
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 31949
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 95847
$ws.Range("N71").Value = -103959
$ws.Range("M71").ClearContents()
# Row 95
$ws.Range("H95").Value = 42124.5
$ws.Range("J95").Value = 42124.5
$ws.Range("L95").Value = 42124.5
$ws.Range("N95").Value = -47616.5
# Row 107
$ws.Range("H107").Value = 613.65
$ws.Range("I107").Value = 600
$ws.Range("J107").Value = 668.25
$ws.Range("K107").Value = 1800
$ws.Range("L107").Value = 2004.75
$ws.Range("M107").Value = 120
$ws.Range("N107").Value = -5844.75
# Row 140
$ws.Range("H140").Value = 43330
$ws.Range("J140").Value = 43330
$ws.Range("L140").Value = 43330
$ws.Range("N140").Value = -53690
# Row 141
$ws.Range("H141").Value = 154432.75
$ws.Range("J141").Value = 154432.75
$ws.Range("L141").Value = 154432.75
$ws.Range("N141").Value = -164792.75

